# Fixed naive component forecaster bug - Presentation state 11.02.
# Rewrites the naive QoQ error series (staircase matrix, rows 2-16 / cols B-K)
# with the corrected forecaster output values, including the newly
# populated diagonal cells that extend each row by one more quarter.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.6713308423207838
$ws.Cells.Item(2, 3).Value = 0.4081742755916882
$ws.Cells.Item(2, 4).Value = 0.18001876026763
$ws.Cells.Item(2, 5).Value = -0.04386758844314925
$ws.Cells.Item(2, 6).Value = 0.8714298707313746
$ws.Cells.Item(2, 7).Value = 0.7287129257210216
$ws.Cells.Item(2, 8).Value = 0.3687093173527261
$ws.Cells.Item(2, 9).Value = 0.5318167280777406
$ws.Cells.Item(2, 10).Value = 0.6850938218533179
$ws.Cells.Item(2, 11).Value = 0.1817452544735015
$ws.Cells.Item(3, 2).Value = 0.4306022552246968
$ws.Cells.Item(3, 3).Value = 0.1826222998300787
$ws.Cells.Item(3, 4).Value = -0.06593527409716449
$ws.Cells.Item(3, 5).Value = 0.8654478904270196
$ws.Cells.Item(3, 6).Value = 0.7233696755386105
$ws.Cells.Item(3, 7).Value = 0.3583407839823598
$ws.Cells.Item(3, 8).Value = 0.5233394059541168
$ws.Cells.Item(3, 9).Value = 0.6773215135465072
$ws.Cells.Item(3, 10).Value = 0.1731212809405329
$ws.Cells.Item(3, 11).Value = 0.4559179537672726
$ws.Cells.Item(4, 2).Value = 0.2040468712872988
$ws.Cells.Item(4, 3).Value = 0.05424257470204791
$ws.Cells.Item(4, 4).Value = 0.7762535460132365
$ws.Cells.Item(4, 5).Value = 0.6936426716900327
$ws.Cells.Item(4, 6).Value = 0.3669669040122699
$ws.Cells.Item(4, 7).Value = 0.4973992325286816
$ws.Cells.Item(4, 8).Value = 0.6545429505182796
$ws.Cells.Item(4, 9).Value = 0.1595032457919151
$ws.Cells.Item(4, 10).Value = 0.4373971565820968
$ws.Cells.Item(4, 11).Value = 0.2648072220727796
$ws.Cells.Item(5, 2).Value = 0.01192194238165845
$ws.Cells.Item(5, 3).Value = 0.7470490721575034
$ws.Cells.Item(5, 4).Value = 0.6906220718851115
$ws.Cells.Item(5, 5).Value = 0.3481056204007895
$ws.Cells.Item(5, 6).Value = 0.4770681014634773
$ws.Cells.Item(5, 7).Value = 0.6400127950840317
$ws.Cells.Item(5, 8).Value = 0.1428622002996471
$ws.Cells.Item(5, 9).Value = 0.4198133111533466
$ws.Cells.Item(5, 10).Value = 0.2483333385427368
$ws.Cells.Item(5, 11).Value = 0.5550469433309027
$ws.Cells.Item(6, 2).Value = 1.087246812775413
$ws.Cells.Item(6, 3).Value = 0.7659931336690486
$ws.Cells.Item(6, 4).Value = 0.1570123340877904
$ws.Cells.Item(6, 5).Value = 0.5017620140502106
$ws.Cells.Item(6, 6).Value = 0.6493921986403277
$ws.Cells.Item(6, 7).Value = 0.08805140168825606
$ws.Cells.Item(6, 8).Value = 0.4009192202646768
$ws.Cells.Item(6, 9).Value = 0.2334693377276911
$ws.Cells.Item(6, 10).Value = 0.526335232701532
$ws.Cells.Item(6, 11).Value = -0.08123716056912761
$ws.Cells.Item(7, 2).Value = 1.216688005659299
$ws.Cells.Item(7, 3).Value = 0.2041230363001488
$ws.Cells.Item(7, 4).Value = 0.2614787832594023
$ws.Cells.Item(7, 5).Value = 0.6853372797061905
$ws.Cells.Item(7, 6).Value = 0.08600064424866319
$ws.Cells.Item(7, 7).Value = 0.3243509808286466
$ws.Cells.Item(7, 8).Value = 0.2060699433524077
$ws.Cells.Item(7, 9).Value = 0.5000461825038066
$ws.Cells.Item(7, 10).Value = -0.1246344512824602
$ws.Cells.Item(7, 11).Value = 0.5737435035592049
$ws.Cells.Item(8, 2).Value = 0.5164486232236872
$ws.Cells.Item(8, 3).Value = 0.3947675229949266
$ws.Cells.Item(8, 4).Value = 0.506258857889999
$ws.Cells.Item(8, 5).Value = 0.1141948830192304
$ws.Cells.Item(8, 6).Value = 0.3603762920210401
$ws.Cells.Item(8, 7).Value = 0.1685283731481941
$ws.Cells.Item(8, 8).Value = 0.4945038434164454
$ws.Cells.Item(8, 9).Value = -0.1204591223437116
$ws.Cells.Item(8, 10).Value = 0.562558789819434
$ws.Cells.Item(9, 2).Value = 0.6303580141027678
$ws.Cells.Item(9, 3).Value = 0.5909375987643086
$ws.Cells.Item(9, 4).Value = -0.03140548361448672
$ws.Cells.Item(9, 5).Value = 0.3703300348802827
$ws.Cells.Item(9, 6).Value = 0.183893542847739
$ws.Cells.Item(9, 7).Value = 0.4556555084590223
$ws.Cells.Item(9, 8).Value = -0.1363049313708975
$ws.Cells.Item(9, 9).Value = 0.5542245613519331
$ws.Cells.Item(10, 2).Value = 0.9019617852456914
$ws.Cells.Item(10, 3).Value = 0.08568329079670847
$ws.Cells.Item(10, 4).Value = 0.2078658934307159
$ws.Cells.Item(10, 5).Value = 0.2128524933129406
$ws.Cells.Item(10, 6).Value = 0.4913983809139259
$ws.Cells.Item(10, 7).Value = -0.1676325777545246
$ws.Cells.Item(10, 8).Value = 0.5509995116504074
$ws.Cells.Item(11, 2).Value = 0.3325084682008229
$ws.Cells.Item(11, 3).Value = 0.2255217609686054
$ws.Cells.Item(11, 4).Value = 0.1182302542019461
$ws.Cells.Item(11, 5).Value = 0.5236157691624059
$ws.Cells.Item(11, 6).Value = -0.1552444519268073
$ws.Cells.Item(11, 7).Value = 0.5301357081011285
$ws.Cells.Item(12, 2).Value = 0.4651981203848173
$ws.Cells.Item(12, 3).Value = 0.2032544976711596
$ws.Cells.Item(12, 4).Value = 0.4071311908043919
$ws.Cells.Item(12, 5).Value = -0.1401391802749548
$ws.Cells.Item(12, 6).Value = 0.5533747437271186
$ws.Cells.Item(13, 2).Value = 0.3681145747052469
$ws.Cells.Item(13, 3).Value = 0.4206684630523081
$ws.Cells.Item(13, 4).Value = -0.2031120992649178
$ws.Cells.Item(13, 5).Value = 0.5634129094086165
$ws.Cells.Item(14, 2).Value = 0.6745175049177161
$ws.Cells.Item(14, 3).Value = -0.1037940490828814
$ws.Cells.Item(14, 4).Value = 0.4518363671933744
$ws.Cells.Item(15, 2).Value = -0.0597343578434324
$ws.Cells.Item(15, 3).Value = 0.4732568720679752
$ws.Cells.Item(16, 2).Value = 0.7095000033804217
